# Auto-update predictions and index for 2025-10-16
# Shifts the fixtures table down by one data row (a new fixture is
# inserted at the top, row 2) and recomputes the tally block in
# columns K/L accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. New data for rows 2-9 (columns A-I). Row 2 is the newly added
#    fixture; rows 3-9 are the previous rows 2-8 shifted down by one,
#    each with refreshed Win% / Result / Total / Under values.
# ---------------------------------------------------------------

$rows = @(
    @{ Row = 2;  A = "Wed Oct 15"; B = "Jamaica ✓ - Bermuda: 4:0";                                   C = 3.62; D = "Jamaica";                        E = 4.5; F = "71%"; G = "✓";  H = 4;  I = $true }
    @{ Row = 3;  A = "Wed Oct 15"; B = "FBC Melgar ✓ - Alianza Universidad: 2:1";                    C = 1.74; D = "FBC Melgar";                     E = 2.5; F = "70%"; G = "✓";  H = 3;  I = $false }
    @{ Row = 4;  A = "Wed Oct 15"; B = "Puerto Rico - Argentina ✓: 0:6";                             C = 1.49; D = "Argentina";                      E = 2.5; F = "68%"; G = "✓";  H = 6;  I = $false }
    @{ Row = 5;  A = "Wed Oct 15"; B = "Sociedade Esportiva Palmeiras ✓ - Red Bull Bragantino: 5:1";  C = 2.04; D = "Sociedade Esportiva Palmeiras"; E = 3.5; F = "68%"; G = "✓";  H = 6;  I = $false }
    @{ Row = 6;  A = "Wed Oct 15"; B = "Club Athletico Paranaense  - Avaí FC: 1:1";                   C = 1.81; D = "Club Athletico Paranaense";      E = 2.5; F = "65%"; G = $null; H = 2;  I = $true }
    @{ Row = 7;  A = "Wed Oct 15"; B = "Al-Nahda Club  - Sur SC: 14:10";                              C = 1;    D = "Al-Nahda Club";                  E = 1.5; F = "61%"; G = $null; H = 24; I = $false }
    @{ Row = 8;  A = "Wed Oct 15"; B = "Atlético Nacional ✓ - Deportivo Cali: 2:1";                   C = 1.36; D = "Atlético Nacional";              E = 2.5; F = "61%"; G = "✓";  H = 3;  I = $false }
    @{ Row = 9;  A = "Wed Oct 15"; B = "United States ✓ - Australia: 2:1";                            C = 1.34; D = "United States";                  E = 2.5; F = "58%"; G = "✓";  H = 3;  I = $false }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    # Leading apostrophe forces text storage so "71%" etc. stay literal
    # strings (matching the source file) instead of turning into a
    # numeric percentage value/format. Resetting the style afterwards
    # drops the "quote prefix" marker Excel would otherwise add, so the
    # cell keeps the plain/default (un-styled) look of the source file.
    $ws.Cells.Item($r.Row, 6).Value = "'" + $r.F
    $ws.Cells.Item($r.Row, 6).Style = "Normal"
    if ($null -eq $r.G) {
        $ws.Cells.Item($r.Row, 7).ClearContents()
    } else {
        $ws.Cells.Item($r.Row, 7).Value = $r.G
    }
    $ws.Cells.Item($r.Row, 8).Value = $r.H
    $ws.Cells.Item($r.Row, 9).Value = $r.I
}

# ---------------------------------------------------------------
# 2. The tally block (previously rows 9-11) moves down to rows
#    10-12, with formulas re-pointed at the new row numbers. Row 9's
#    old K/L formulas are cleared since row 9 is now a data row.
# ---------------------------------------------------------------

$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()

$ws.Range("K10").Formula = "=COUNTIF(I:I,TRUE)"
$ws.Range("L10").Formula = "=(K10/K12)*100"
$ws.Range("K11").Formula = "=COUNTIF(I:I,FALSE)"
$ws.Range("K12").Formula = "=K10+K11"

# ---------------------------------------------------------------
# 3. Update the sheet's used-range dimension to match the extra row.
# ---------------------------------------------------------------

$ws.UsedRange | Out-Null
